$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in R4 grades (column E) for all students
$ws.Range("E2").Value = 0
$ws.Range("E3").Value = 0
$ws.Range("E4").Value = 2.5
$ws.Range("E5").Value = 0
$ws.Range("E6").Value = 2.5

# Fill in Conceito (column G) for all students
$ws.Range("G2").Value = "RF"
$ws.Range("G3").Value = "A"
$ws.Range("G4").Value = "A"
$ws.Range("G5").Value = "RF"
$ws.Range("G6").Value = "A"

# Update selection to match final cursor position
$ws.Range("G7").Select()
